# PGx review update for sample 20220419-630044
#
# The marker CYP2D6_4 / CYP2D6_011 (sample S2) allele "A" peak was
# re-reviewed by the PGx team: its minimum-bin/height thresholds were
# widened, which lets the previously-missed "A" peak be detected. This
# flips the call for that marker from homozygous wildtype (G/G) to
# heterozygous (G/A), and the overall sample genotype from *1/*10B to
# *1/*4.
#
# This script reproduces the underlying data edits across the four
# linked result sheets (peak_table, allele_table, marker_table,
# genotype_result) plus a couple of harmless workbook-level touch-ups
# (dropping the stale workbookProtection marker and re-selecting the
# cell that was actively being edited).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# peak_table: widen the w_min bin and correct the reported peak heights
# for CYP2D6_4 (row 12).
# ---------------------------------------------------------------------
$peak = $wb.Worksheets.Item("peak_table")
$peak.Range("F12").Value = 27
$peak.Range("N12").Value = 500
$peak.Range("O12").Value = 700

# ---------------------------------------------------------------------
# allele_table: same marker, two allele rows (G=wildtype call, A=mutant
# call). Row 22 = "G" allele, row 23 = "A" allele.
# ---------------------------------------------------------------------
$allele = $wb.Worksheets.Item("allele_table")

# Row 22 ("G" allele) - bin/height numbers updated, peak now resolved.
$allele.Range("I22").Value = 27
$allele.Range("K22").Value = 500
$allele.Range("N22").Value = 36
$allele.Range("O22").Value = 30.96
$allele.Range("P22").Value = 558

# Row 23 ("A" allele) - previously undetected, now detected as a real
# peak with the PGx team's revised thresholds.
$allele.Range("K23").Value = 700
$allele.Range("M23").Value = $true
$allele.Range("N23").Value = 31
$allele.Range("O23").Value = 32.35
$allele.Range("P23").Value = 751
$allele.Range("Q23").Value = "ok"
$allele.Range("R23").ClearContents()

# ---------------------------------------------------------------------
# marker_table: genotype call for CYP2D6_4 flips from homozygous
# wildtype to heterozygous now that the "A" allele is detected.
# ---------------------------------------------------------------------
$marker = $wb.Worksheets.Item("marker_table")
$marker.Range("G12").Value = "GA"
$marker.Range("H12").Value = "heterozygous"

# ---------------------------------------------------------------------
# genotype_result: overall sample genotype updated accordingly.
# ---------------------------------------------------------------------
$genotype = $wb.Worksheets.Item("genotype_result")
$genotype.Range("B2").Value = "*1/*4"

# ---------------------------------------------------------------------
# Workbook-level touch-ups mirroring the re-save: the workbook was never
# actually protected, so drop the now-stale protection marker, and leave
# the selection on the cell that was last edited during the review.
# ---------------------------------------------------------------------
$wb.Unprotect()

$peak.Activate()
$peak.Range("N12").Select()
